# Auto-generated Excel COM-interop script to apply the scraper-refresh diff
# to '苏州-漫展信息.xlsx' (Suzhou comic-convention info workbook).
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 (Exhibitions) ----
$ws1 = $wb.Worksheets.Item(1)

# Bump the live want-to-go counter for row 2
$ws1.Range("F2").Value2 = 796

# Insert a newly-scraped event as row 5, pushing the rest down by one
$ws1.Rows.Item(5).Insert()

# Copy formatting of the index cell (column A) down from the row above
$ws1.Range("A4").Copy()
$ws1.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 5
$ws1.Range("A5").Value2 = 4
$ws1.Range("B5").Value2 = '2024-10-19'
$ws1.Range("C5").Value2 = '苏州·创世次元动漫游戏嘉年华X嘉宾大咖签售会X快来加入这场狂欢，一起嗨翻全场！（免费展）'
$ws1.Range("D5").Value2 = '中山南路1818号 苏州吴江天虹购物中心'
$ws1.Range("E5").Value2 = '2024.10.19 10:30-10.20 20:00'
$ws1.Range("F5").Value2 = 5
$ws1.Range("G5").Value2 = 58
$ws1.Range("H5").Value2 = 'https://show.bilibili.com/platform/detail.html?id=93201'
$ws1.Range("I5").Value2 = '//i1.hdslb.com/bfs/openplatform/202410/WSDkTnxT1728455264113.jpeg'

# Row 6
$ws1.Range("A6").Value2 = 5
$ws1.Range("B6").Value2 = '2024-10-19'
$ws1.Range("C6").Value2 = '苏州·无限次元夜场'
$ws1.Range("D6").Value2 = '盘胥路826号 sos时尚国际CLUB'
$ws1.Range("E6").Value2 = '2024.10.19 13:00-10.19 19:00'
$ws1.Range("F6").Value2 = 157
$ws1.Range("G6").Value2 = 169
$ws1.Range("H6").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92391'
$ws1.Range("I6").Value2 = '//i2.hdslb.com/bfs/openplatform/202409/QaX6Akvh1726131978180.png'

# Row 7
$ws1.Range("A7").Value2 = 6
$ws1.Range("B7").Value2 = '2024-10-26'
$ws1.Range("C7").Value2 = '苏州·国内知名配音演员吕书君@阿君归来 签售会'
$ws1.Range("D7").Value2 = '开平路2188号 苏州吾悦广场'
$ws1.Range("E7").Value2 = '2024.10.26 13:30-10.26 16:00'
$ws1.Range("F7").Value2 = 20
$ws1.Range("G7").Value2 = 125
$ws1.Range("H7").Value2 = 'https://show.bilibili.com/platform/detail.html?id=93060'
$ws1.Range("I7").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/ycNGgQIK1727504415140.png'

# Row 8
$ws1.Range("A8").Value2 = 7
$ws1.Range("B8").Value2 = '2024-10-26'
$ws1.Range("C8").Value2 = '苏州·漫语堂动漫嘉年华'
$ws1.Range("D8").Value2 = '金山东路78号 苏州狮山国际会展中心'
$ws1.Range("E8").Value2 = '2024.10.26 10:00-10.27 17:00'
$ws1.Range("F8").Value2 = 180
$ws1.Range("G8").Value2 = 70
$ws1.Range("H8").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91626'
$ws1.Range("I8").Value2 = '//i1.hdslb.com/bfs/openplatform/202408/HxlG5vVw1724918564912.jpeg'

# Row 9
$ws1.Range("A9").Value2 = 8
$ws1.Range("B9").Value2 = '2024-10-26'
$ws1.Range("C9").Value2 = '苏州·第三届华盟国漫次元嘉年华'
$ws1.Range("D9").Value2 = '清禾路886号 苏州聚橙尹山湖大剧院'
$ws1.Range("E9").Value2 = '2024.10.26 10:00-10.27 17:00'
$ws1.Range("F9").Value2 = 365
$ws1.Range("G9").Value2 = 58
$ws1.Range("H9").Value2 = 'https://show.bilibili.com/platform/detail.html?id=85767'
$ws1.Range("I9").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/52AMZyUi1727059410434.jpeg'

# Row 10
$ws1.Range("A10").Value2 = 9
$ws1.Range("B10").Value2 = '2024-11-02'
$ws1.Range("C10").Value2 = '苏州·女神异闻录only同人展'
$ws1.Range("D10").Value2 = '小外滩街苏州小外滩 元和塘美术馆'
$ws1.Range("E10").Value2 = '2024.11.02 09:00-11.03 17:00'
$ws1.Range("F10").Value2 = 476
$ws1.Range("G10").Value2 = 78
$ws1.Range("H10").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91709'
$ws1.Range("I10").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/PpYo3LGQ1725214760478.jpeg'

# Row 11
$ws1.Range("A11").Value2 = 10
$ws1.Range("B11").Value2 = '2024-11-02'
$ws1.Range("C11").Value2 = '苏州·绘时国乙2.0光夜同人only'
$ws1.Range("D11").Value2 = '东亭街588号 南舍别院'
$ws1.Range("E11").Value2 = '2024.11.02 10:30-11.02 20:30'
$ws1.Range("F11").Value2 = 520
$ws1.Range("G11").Value2 = 178
$ws1.Range("H11").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91324'
$ws1.Range("I11").Value2 = '//i1.hdslb.com/bfs/openplatform/202408/YauAhbAd1724662566605.jpeg'

# Row 12
$ws1.Range("A12").Value2 = 11
$ws1.Range("B12").Value2 = '2024-11-16'
$ws1.Range("C12").Value2 = '张家港·META萌圆饿了'
$ws1.Range("D12").Value2 = '杨舍镇人民中路42号 张家港国贸酒店'
$ws1.Range("E12").Value2 = '2024.11.16 10:00-11.16 17:00'
$ws1.Range("F12").Value2 = 148
$ws1.Range("G12").Value2 = 40
$ws1.Range("H12").Value2 = 'https://show.bilibili.com/platform/detail.html?id=90745'
$ws1.Range("I12").Value2 = '//i2.hdslb.com/bfs/openplatform/202408/jB7b2kZ11723621730632.jpeg'

# Row 13
$ws1.Range("A13").Value2 = 12
$ws1.Range("B13").Value2 = '2024-11-16'
$ws1.Range("C13").Value2 = '苏州·COME IN JOY 动漫品牌国潮文化节'
$ws1.Range("D13").Value2 = '金山南路288号 木渎影视城会展中心'
$ws1.Range("E13").Value2 = '2024.11.16 10:00-11.17 17:00'
$ws1.Range("F13").Value2 = 12108
$ws1.Range("G13").Value2 = 60
$ws1.Range("H13").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92177'
$ws1.Range("I13").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/pBNpUxOr1725866134021.jpeg'

# Row 14
$ws1.Range("A14").Value2 = 13
$ws1.Range("B14").Value2 = '2025-01-01'
$ws1.Range("C14").Value2 = '苏州·星部落&青铜树动漫嘉年华'
$ws1.Range("D14").Value2 = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Range("E14").Value2 = '2025.01.01 09:00-01.02 16:00'
$ws1.Range("F14").Value2 = 5445
$ws1.Range("G14").Value2 = 55
$ws1.Range("H14").Value2 = 'https://show.bilibili.com/platform/detail.html?id=84858'
$ws1.Range("I14").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/u3RjLCRL1727662424227.jpeg'

# ---- Sheet 2: 演出 (Performances) ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value2 = 109

# ---- Sheet 4: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value2 = 796
$ws4.Range("F4").Value2 = 109

# Insert the same newly-scraped event as row 7, pushing the rest down by one
$ws4.Rows.Item(7).Insert()

$ws4.Range("A6").Copy()
$ws4.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 7
$ws4.Range("A7").Value2 = 6
$ws4.Range("B7").Value2 = '2024-10-19'
$ws4.Range("C7").Value2 = '苏州·创世次元动漫游戏嘉年华X嘉宾大咖签售会X快来加入这场狂欢，一起嗨翻全场！（免费展）'
$ws4.Range("D7").Value2 = '中山南路1818号 苏州吴江天虹购物中心'
$ws4.Range("E7").Value2 = '2024.10.19 10:30-10.20 20:00'
$ws4.Range("F7").Value2 = 5
$ws4.Range("G7").Value2 = 58
$ws4.Range("H7").Value2 = 'https://show.bilibili.com/platform/detail.html?id=93201'
$ws4.Range("I7").Value2 = '//i1.hdslb.com/bfs/openplatform/202410/WSDkTnxT1728455264113.jpeg'

# Row 8
$ws4.Range("A8").Value2 = 7
$ws4.Range("B8").Value2 = '2024-10-19'
$ws4.Range("C8").Value2 = '苏州·无限次元夜场'
$ws4.Range("D8").Value2 = '盘胥路826号 sos时尚国际CLUB'
$ws4.Range("E8").Value2 = '2024.10.19 13:00-10.19 19:00'
$ws4.Range("F8").Value2 = 157
$ws4.Range("G8").Value2 = 169
$ws4.Range("H8").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92391'
$ws4.Range("I8").Value2 = '//i2.hdslb.com/bfs/openplatform/202409/QaX6Akvh1726131978180.png'

# Row 9
$ws4.Range("A9").Value2 = 8
$ws4.Range("B9").Value2 = '2024-10-26'
$ws4.Range("C9").Value2 = '苏州·国内知名配音演员吕书君@阿君归来 签售会'
$ws4.Range("D9").Value2 = '开平路2188号 苏州吾悦广场'
$ws4.Range("E9").Value2 = '2024.10.26 13:30-10.26 16:00'
$ws4.Range("F9").Value2 = 20
$ws4.Range("G9").Value2 = 125
$ws4.Range("H9").Value2 = 'https://show.bilibili.com/platform/detail.html?id=93060'
$ws4.Range("I9").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/ycNGgQIK1727504415140.png'

# Row 10
$ws4.Range("A10").Value2 = 9
$ws4.Range("B10").Value2 = '2024-10-26'
$ws4.Range("C10").Value2 = '苏州·漫语堂动漫嘉年华'
$ws4.Range("D10").Value2 = '金山东路78号 苏州狮山国际会展中心'
$ws4.Range("E10").Value2 = '2024.10.26 10:00-10.27 17:00'
$ws4.Range("F10").Value2 = 180
$ws4.Range("G10").Value2 = 70
$ws4.Range("H10").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91626'
$ws4.Range("I10").Value2 = '//i1.hdslb.com/bfs/openplatform/202408/HxlG5vVw1724918564912.jpeg'

# Row 11
$ws4.Range("A11").Value2 = 10
$ws4.Range("B11").Value2 = '2024-10-26'
$ws4.Range("C11").Value2 = '苏州·第三届华盟国漫次元嘉年华'
$ws4.Range("D11").Value2 = '清禾路886号 苏州聚橙尹山湖大剧院'
$ws4.Range("E11").Value2 = '2024.10.26 10:00-10.27 17:00'
$ws4.Range("F11").Value2 = 365
$ws4.Range("G11").Value2 = 58
$ws4.Range("H11").Value2 = 'https://show.bilibili.com/platform/detail.html?id=85767'
$ws4.Range("I11").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/52AMZyUi1727059410434.jpeg'

# Row 12
$ws4.Range("A12").Value2 = 11
$ws4.Range("B12").Value2 = '2024-11-02'
$ws4.Range("C12").Value2 = '苏州·女神异闻录only同人展'
$ws4.Range("D12").Value2 = '小外滩街苏州小外滩 元和塘美术馆'
$ws4.Range("E12").Value2 = '2024.11.02 09:00-11.03 17:00'
$ws4.Range("F12").Value2 = 476
$ws4.Range("G12").Value2 = 78
$ws4.Range("H12").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91709'
$ws4.Range("I12").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/PpYo3LGQ1725214760478.jpeg'

# Row 13
$ws4.Range("A13").Value2 = 12
$ws4.Range("B13").Value2 = '2024-11-02'
$ws4.Range("C13").Value2 = '苏州·绘时国乙2.0光夜同人only'
$ws4.Range("D13").Value2 = '东亭街588号 南舍别院'
$ws4.Range("E13").Value2 = '2024.11.02 10:30-11.02 20:30'
$ws4.Range("F13").Value2 = 520
$ws4.Range("G13").Value2 = 178
$ws4.Range("H13").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91324'
$ws4.Range("I13").Value2 = '//i1.hdslb.com/bfs/openplatform/202408/YauAhbAd1724662566605.jpeg'

# Row 14
$ws4.Range("A14").Value2 = 13
$ws4.Range("B14").Value2 = '2024-11-16'
$ws4.Range("C14").Value2 = '张家港·META萌圆饿了'
$ws4.Range("D14").Value2 = '杨舍镇人民中路42号 张家港国贸酒店'
$ws4.Range("E14").Value2 = '2024.11.16 10:00-11.16 17:00'
$ws4.Range("F14").Value2 = 148
$ws4.Range("G14").Value2 = 40
$ws4.Range("H14").Value2 = 'https://show.bilibili.com/platform/detail.html?id=90745'
$ws4.Range("I14").Value2 = '//i2.hdslb.com/bfs/openplatform/202408/jB7b2kZ11723621730632.jpeg'

# Row 15
$ws4.Range("A15").Value2 = 14
$ws4.Range("B15").Value2 = '2024-11-16'
$ws4.Range("C15").Value2 = '苏州·COME IN JOY 动漫品牌国潮文化节'
$ws4.Range("D15").Value2 = '金山南路288号 木渎影视城会展中心'
$ws4.Range("E15").Value2 = '2024.11.16 10:00-11.17 17:00'
$ws4.Range("F15").Value2 = 12108
$ws4.Range("G15").Value2 = 60
$ws4.Range("H15").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92177'
$ws4.Range("I15").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/pBNpUxOr1725866134021.jpeg'

# Row 16
$ws4.Range("A16").Value2 = 15
$ws4.Range("B16").Value2 = '2024-12-22'
$ws4.Range("C16").Value2 = '苏州·维也纳皇家交响乐团2025新年音乐会'
$ws4.Range("D16").Value2 = '东苑路1号公共文化中心内 苏州保利大剧院'
$ws4.Range("E16").Value2 = '2024.12.22 19:30-12.22 21:30'
$ws4.Range("F16").Value2 = 8
$ws4.Range("G16").Value2 = 280
$ws4.Range("H16").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92817'
$ws4.Range("I16").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/cCTiFEpg1727155421223.jpeg'

# Row 17
$ws4.Range("A17").Value2 = 16
$ws4.Range("B17").Value2 = '2025-01-01'
$ws4.Range("C17").Value2 = '苏州·星部落&青铜树动漫嘉年华'
$ws4.Range("D17").Value2 = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Range("E17").Value2 = '2025.01.01 09:00-01.02 16:00'
$ws4.Range("F17").Value2 = 5445
$ws4.Range("G17").Value2 = 55
$ws4.Range("H17").Value2 = 'https://show.bilibili.com/platform/detail.html?id=84858'
$ws4.Range("I17").Value2 = '//i1.hdslb.com/bfs/openplatform/202409/u3RjLCRL1727662424227.jpeg'

Write-Host "Edit complete"